# Insert a new weekly record row above row 50 ("Fruta / hortaliza, semanal"),
# shifting all subsequent rows down by one (old row 50 -> new row 51, ...,
# old row 124 -> new row 125).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 50; existing rows 50-124 shift to 51-125.
$ws.Rows.Item(50).Insert()

# Populate the newly inserted row 50 with the new weekly record.
$ws.Cells.Item(50, 1).Value2 = 6
$ws.Cells.Item(50, 2).Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(50, 3).Value2 = "Metropolitana"
$ws.Cells.Item(50, 4).Value2 = 44949
$ws.Cells.Item(50, 5).Value2 = 13
$ws.Cells.Item(50, 6).Value2 = "Fruta"
$ws.Cells.Item(50, 7).Value2 = 100101
$ws.Cells.Item(50, 8).Value2 = "Berries"
$ws.Cells.Item(50, 9).Value2 = 100101008
$ws.Cells.Item(50, 10).Value2 = "Mora"
$ws.Cells.Item(50, 11).Value2 = "Sin especificar"
$ws.Cells.Item(50, 12).Value2 = "Primera"
$ws.Cells.Item(50, 13).Value2 = 150
$ws.Cells.Item(50, 14).Value2 = 4000
$ws.Cells.Item(50, 15).Value2 = 4000
$ws.Cells.Item(50, 16).Value2 = 4000
$ws.Cells.Item(50, 17).Value2 = "`$/bandeja 2 kilos"
$ws.Cells.Item(50, 18).Value2 = "Región del Maule"
$ws.Cells.Item(50, 19).Value2 = 2000
$ws.Cells.Item(50, 20).Value2 = 2
